{"js": "// Replace the date line and the 25 multiplication problems with their new values.\n// All source strings are unique in the document, so a search+replace per pair is safe.\nconst replacements = [\n  [\"2024-12-16 Monday\", \"2024-12-17 Tuesday\"],\n  [\"444\u00d74=\", \"376\u00d73=\"],\n  [\"994\u00d75=\", \"429\u00d72=\"],\n  [\"713\u00d72=\", \"840\u00d77=\"],\n  [\"788\u00d75=\", \"907\u00d74=\"],\n  [\"779\u00d75=\", \"246\u00d76=\"],\n  [\"395\u00d79=\", \"603\u00d75=\"],\n  [\"896\u00d79=\", \"466\u00d78=\"],\n  [\"677\u00d79=\", \"559\u00d75=\"],\n  [\"843\u00d74=\", \"666\u00d72=\"],\n  [\"755\u00d76=\", \"132\u00d73=\"],\n  [\"456\u00d77=\", \"978\u00d76=\"],\n  [\"937\u00d77=\", \"845\u00d77=\"],\n  [\"587\u00d77=\", \"214\u00d75=\"],\n  [\"153\u00d77=\", \"271\u00d78=\"],\n  [\"428\u00d75=\", \"855\u00d77=\"],\n  [\"664\u00d72=\", \"989\u00d76=\"],\n  [\"956\u00d77=\", \"461\u00d77=\"],\n  [\"997\u00d72=\", \"474\u00d77=\"],\n  [\"910\u00d77=\", \"959\u00d78=\"],\n  [\"900\u00d77=\", \"389\u00d79=\"],\n  [\"854\u00d79=\", \"822\u00d78=\"],\n  [\"471\u00d78=\", \"530\u00d73=\"],\n  [\"350\u00d75=\", \"840\u00d76=\"],\n  [\"759\u00d74=\", \"223\u00d73=\"],\n  [\"489\u00d75=\", \"941\u00d77=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2024-12-16 Monday\", \"2024-12-17 Tuesday\"),\n    @(\"444\u00d74=\", \"376\u00d73=\"),\n    @(\"994\u00d75=\", \"429\u00d72=\"),\n    @(\"713\u00d72=\", \"840\u00d77=\"),\n    @(\"788\u00d75=\", \"907\u00d74=\"),\n    @(\"779\u00d75=\", \"246\u00d76=\"),\n    @(\"395\u00d79=\", \"603\u00d75=\"),\n    @(\"896\u00d79=\", \"466\u00d78=\"),\n    @(\"677\u00d79=\", \"559\u00d75=\"),\n    @(\"843\u00d74=\", \"666\u00d72=\"),\n    @(\"755\u00d76=\", \"132\u00d73=\"),\n    @(\"456\u00d77=\", \"978\u00d76=\"),\n    @(\"937\u00d77=\", \"845\u00d77=\"),\n    @(\"587\u00d77=\", \"214\u00d75=\"),\n    @(\"153\u00d77=\", \"271\u00d78=\"),\n    @(\"428\u00d75=\", \"855\u00d77=\"),\n    @(\"664\u00d72=\", \"989\u00d76=\"),\n    @(\"956\u00d77=\", \"461\u00d77=\"),\n    @(\"997\u00d72=\", \"474\u00d77=\"),\n    @(\"910\u00d77=\", \"959\u00d78=\"),\n    @(\"900\u00d77=\", \"389\u00d79=\"),\n    @(\"854\u00d79=\", \"822\u00d78=\"),\n    @(\"471\u00d78=\", \"530\u00d73=\"),\n    @(\"350\u00d75=\", \"840\u00d76=\"),\n    @(\"759\u00d74=\", \"223\u00d73=\"),\n    @(\"489\u00d75=\", \"941\u00d77=\")\n)\n\nforeach ($pair in $replacements) {\n    $old = $pair[0]\n    $new = $pair[1]\n\n    $rng = $d.Content\n    $find = $rng.Find\n    $find.ClearFormatting()\n    $find.Text = $old\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $new\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.Execute([ref]$old, [ref]$true, [ref]$false, [ref]$false, [ref]$false, [ref]$false, [ref]$true, 1, [ref]$false, $new, 2)\n}\n"}
